{"js": "// The commit swaps every two-digit-by-two-digit multiplication expression\n// in the practice-sheet table for a freshly generated one. Each \"old\" string\n// below occurs exactly once in the document, so a literal (non-wildcard)\n// search-and-replace per pair reproduces the diff unambiguously.\nconst pairs = [\n  [\"19\u00d777=1463\", \"28\u00d772=2016\"],\n  [\"67\u00d772=4824\", \"95\u00d719=1805\"],\n  [\"94\u00d724=2256\", \"64\u00d719=1216\"],\n  [\"95\u00d784=7980\", \"42\u00d776=3192\"],\n  [\"83\u00d794=7802\", \"38\u00d711=418\"],\n  [\"13\u00d783=1079\", \"39\u00d754=2106\"],\n  [\"92\u00d795=8740\", \"77\u00d781=6237\"],\n  [\"12\u00d792=1104\", \"72\u00d782=5904\"],\n  [\"22\u00d787=1914\", \"17\u00d749=833\"],\n  [\"25\u00d759=1475\", \"66\u00d761=4026\"],\n  [\"50\u00d767=3350\", \"34\u00d715=510\"],\n  [\"44\u00d799=4356\", \"34\u00d754=1836\"],\n  [\"64\u00d726=1664\", \"43\u00d771=3053\"],\n  [\"81\u00d750=4050\", \"60\u00d795=5700\"],\n  [\"27\u00d780=2160\", \"75\u00d743=3225\"],\n  [\"13\u00d778=1014\", \"81\u00d743=3483\"],\n  [\"46\u00d759=2714\", \"55\u00d764=3520\"],\n  [\"55\u00d720=1100\", \"79\u00d788=6952\"],\n  [\"52\u00d730=1560\", \"55\u00d771=3905\"],\n  [\"13\u00d765=845\", \"54\u00d784=4536\"],\n  [\"70\u00d745=3150\", \"26\u00d788=2288\"],\n  [\"56\u00d749=2744\", \"65\u00d722=1430\"],\n  [\"20\u00d768=1360\", \"21\u00d754=1134\"],\n  [\"74\u00d715=1110\", \"32\u00d724=768\"],\n  [\"34\u00d766=2244\", \"66\u00d712=792\"]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The commit swaps every two-digit-by-two-digit multiplication expression\n# in the practice-sheet table for a freshly generated one. Each \"old\" string\n# occurs exactly once in the document, so a literal Find/Replace (wdReplaceAll,\n# but effectively \"replace the one match\") per pair reproduces the diff.\n$pairs = @(\n  @(\"19\u00d777=1463\", \"28\u00d772=2016\"),\n  @(\"67\u00d772=4824\", \"95\u00d719=1805\"),\n  @(\"94\u00d724=2256\", \"64\u00d719=1216\"),\n  @(\"95\u00d784=7980\", \"42\u00d776=3192\"),\n  @(\"83\u00d794=7802\", \"38\u00d711=418\"),\n  @(\"13\u00d783=1079\", \"39\u00d754=2106\"),\n  @(\"92\u00d795=8740\", \"77\u00d781=6237\"),\n  @(\"12\u00d792=1104\", \"72\u00d782=5904\"),\n  @(\"22\u00d787=1914\", \"17\u00d749=833\"),\n  @(\"25\u00d759=1475\", \"66\u00d761=4026\"),\n  @(\"50\u00d767=3350\", \"34\u00d715=510\"),\n  @(\"44\u00d799=4356\", \"34\u00d754=1836\"),\n  @(\"64\u00d726=1664\", \"43\u00d771=3053\"),\n  @(\"81\u00d750=4050\", \"60\u00d795=5700\"),\n  @(\"27\u00d780=2160\", \"75\u00d743=3225\"),\n  @(\"13\u00d778=1014\", \"81\u00d743=3483\"),\n  @(\"46\u00d759=2714\", \"55\u00d764=3520\"),\n  @(\"55\u00d720=1100\", \"79\u00d788=6952\"),\n  @(\"52\u00d730=1560\", \"55\u00d771=3905\"),\n  @(\"13\u00d765=845\", \"54\u00d784=4536\"),\n  @(\"70\u00d745=3150\", \"26\u00d788=2288\"),\n  @(\"56\u00d749=2744\", \"65\u00d722=1430\"),\n  @(\"20\u00d768=1360\", \"21\u00d754=1134\"),\n  @(\"74\u00d715=1110\", \"32\u00d724=768\"),\n  @(\"34\u00d766=2244\", \"66\u00d712=792\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n  #   ReplaceWith, Replace=wdReplaceAll)\n  $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
